$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.464.55'
$ws.Range('E2').Value = '  +1.93%  '
$ws.Range('D3').Value = '1.828.81'
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'315.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').Value = "'0.5181"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.43%  '
$ws.Range('D8').Value = "'0.3912"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('D9').Value = "'0.07643"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').Value = "'1.109"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.20%  '
$ws.Range('D12').Value = "'21.12"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.81%  '
$ws.Range('D13').Value = "'6.301"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').Value = "'1.000"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = "'7.552"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.76%  '
$ws.Range('D16').Value = '1.825.62'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').Value = "'93.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.79%  '
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').Value = "'0.06669"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('D20').Value = "'17.71"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('D21').Value = "'1.000"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = "'6.190"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.96%  '
$ws.Range('D23').Value = '28.490.31'
$ws.Range('E23').Value = '  +1.92%  '
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +7.75%  '
$ws.Range('D26').Value = "'156.86"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = "'20.65"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.40%  '
$ws.Range('D28').Value = '2.035.19'
$ws.Range('E28').Value = '  +1.60%  '
$ws.Range('D29').Value = "'2.401"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.70%  '
$ws.Range('D30').Value = "'125.12"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.62%  '
$ws.Range('D31').Value = "'1.124"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.71%  '
$ws.Range('D32').Value = "'0.1087"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('D33').Value = "'5.683"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.29%  '
$ws.Range('D34').Value = "'3.660"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('D35').Value = "'0.07033"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').Value = "'0.2232"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').Value = "'8.978"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.26%  '
$ws.Range('D38').Value = "'0.02327"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').Value = "'5.138"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.67%  '
$ws.Range('D40').Value = "'0.6292"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').Value = "'1.182"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -1.50%  '
$ws.Range('D45').Value = "'13.40"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').Value = "'0.5908"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.49%  '
$ws.Range('D47').Value = "'3.711"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('D48').Value = "'124.32"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').Value = "'1.982"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.39%  '
$ws.Range('D50').Value = "'1.201"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('D51').Value = "'0.06928"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.81%  '
